# Fix tests when C_d is non-zero
#
# The "EEU data" example workbook had its disposal-cost test values
# (C_d_orig in column V, C_d_star in column X, row 2) sitting at 0,
# which doesn't exercise the non-zero C_d code path in the package's
# tests. Populate them with non-zero sample values and leave the
# selection where the author left it (X3) when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C_d_orig (row 2, "Original" case) -> 100
$ws.Range("V2").Value = 100

# C_d_star (row 2, "Original" case) -> 101
$ws.Range("X2").Value = 101

# Match the active-cell selection recorded in the saved workbook.
$ws.Range("X3").Select()
